$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Style the new rows (11-19) to match existing data-row style (border + centered)
$newRows = $ws.Range("A11:J19")
$newRows.HorizontalAlignment = -4108
$newRows.VerticalAlignment = -4108
$newRows.Borders.LineStyle = 1
$newRows.Borders.Weight = 2
$newRows.Borders.Color = 0

# Row 2
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 'Khoa Nông nghiệp'
$ws.Cells.Item(2, 3).Value = 2015
$ws.Cells.Item(2, 4).Value = 'X'
$ws.Cells.Item(2, 5).Value = ''
$ws.Cells.Item(2, 6).Value = ''
$ws.Cells.Item(2, 7).Value = ''
$ws.Cells.Item(2, 8).Value = ''
$ws.Cells.Item(2, 9).Value = ''
$ws.Cells.Item(2, 10).Value = ''

# Row 3
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 'Phòng Kế hoạch - Tài chính'
$ws.Cells.Item(3, 3).Value = 2017
$ws.Cells.Item(3, 4).Value = ''
$ws.Cells.Item(3, 5).Value = 'X'
$ws.Cells.Item(3, 6).Value = ''
$ws.Cells.Item(3, 7).Value = ''
$ws.Cells.Item(3, 8).Value = ''
$ws.Cells.Item(3, 9).Value = ''
$ws.Cells.Item(3, 10).Value = ''

# Row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 'Phòng Quản lý đào tạo'
$ws.Cells.Item(4, 3).Value = 2017
$ws.Cells.Item(4, 4).Value = ''
$ws.Cells.Item(4, 5).Value = 'X'
$ws.Cells.Item(4, 6).Value = ''
$ws.Cells.Item(4, 7).Value = ''
$ws.Cells.Item(4, 8).Value = ''
$ws.Cells.Item(4, 9).Value = ''
$ws.Cells.Item(4, 10).Value = ''

# Row 5
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = 'Phòng Kế hoạch - Tài chính'
$ws.Cells.Item(5, 3).Value = 2017
$ws.Cells.Item(5, 4).Value = ''
$ws.Cells.Item(5, 5).Value = ''
$ws.Cells.Item(5, 6).Value = ''
$ws.Cells.Item(5, 7).Value = 'X'
$ws.Cells.Item(5, 8).Value = ''
$ws.Cells.Item(5, 9).Value = ''
$ws.Cells.Item(5, 10).Value = ''

# Row 6
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 'Phòng Quản lý đào tạo'
$ws.Cells.Item(6, 3).Value = 2017
$ws.Cells.Item(6, 4).Value = ''
$ws.Cells.Item(6, 5).Value = ''
$ws.Cells.Item(6, 6).Value = ''
$ws.Cells.Item(6, 7).Value = 'X'
$ws.Cells.Item(6, 8).Value = ''
$ws.Cells.Item(6, 9).Value = ''
$ws.Cells.Item(6, 10).Value = ''

# Row 7
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = 'Phòng Quản lý Khoa học - Hợp tác'
$ws.Cells.Item(7, 3).Value = 2017
$ws.Cells.Item(7, 4).Value = ''
$ws.Cells.Item(7, 5).Value = ''
$ws.Cells.Item(7, 6).Value = ''
$ws.Cells.Item(7, 7).Value = 'X'
$ws.Cells.Item(7, 8).Value = ''
$ws.Cells.Item(7, 9).Value = ''
$ws.Cells.Item(7, 10).Value = ''

# Row 8
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = 'Trung tâm Kỹ thuật -Nông nghiệp'
$ws.Cells.Item(8, 3).Value = 2019
$ws.Cells.Item(8, 4).Value = ''
$ws.Cells.Item(8, 5).Value = 'X'
$ws.Cells.Item(8, 6).Value = ''
$ws.Cells.Item(8, 7).Value = ''
$ws.Cells.Item(8, 8).Value = ''
$ws.Cells.Item(8, 9).Value = ''
$ws.Cells.Item(8, 10).Value = ''

# Row 9
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 'Trung tâm Ngoại ngữ - Tin học'
$ws.Cells.Item(9, 3).Value = 2019
$ws.Cells.Item(9, 4).Value = ''
$ws.Cells.Item(9, 5).Value = 'X'
$ws.Cells.Item(9, 6).Value = ''
$ws.Cells.Item(9, 7).Value = ''
$ws.Cells.Item(9, 8).Value = ''
$ws.Cells.Item(9, 9).Value = ''
$ws.Cells.Item(9, 10).Value = ''

# Row 10
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = 'Trung tâm Ngoại ngữ - Tin học'
$ws.Cells.Item(10, 3).Value = 2020
$ws.Cells.Item(10, 4).Value = ''
$ws.Cells.Item(10, 5).Value = 'X'
$ws.Cells.Item(10, 6).Value = ''
$ws.Cells.Item(10, 7).Value = ''
$ws.Cells.Item(10, 8).Value = ''
$ws.Cells.Item(10, 9).Value = ''
$ws.Cells.Item(10, 10).Value = ''

# Row 11
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = 'Trung tâm Ngoại ngữ - Tin học'
$ws.Cells.Item(11, 3).Value = '2021-2022'
$ws.Cells.Item(11, 4).Value = ''
$ws.Cells.Item(11, 5).Value = 'X'
$ws.Cells.Item(11, 6).Value = ''
$ws.Cells.Item(11, 7).Value = ''
$ws.Cells.Item(11, 8).Value = ''
$ws.Cells.Item(11, 9).Value = ''
$ws.Cells.Item(11, 10).Value = ''

# Row 12
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 'Trung tâm Ngoại ngữ - Tin học'
$ws.Cells.Item(12, 3).Value = '2022-2023'
$ws.Cells.Item(12, 4).Value = ''
$ws.Cells.Item(12, 5).Value = 'X'
$ws.Cells.Item(12, 6).Value = ''
$ws.Cells.Item(12, 7).Value = ''
$ws.Cells.Item(12, 8).Value = ''
$ws.Cells.Item(12, 9).Value = ''
$ws.Cells.Item(12, 10).Value = ''

# Row 13
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = 'Khoa Công nghệ - Thủy sản'
$ws.Cells.Item(13, 3).Value = '2023-2024'
$ws.Cells.Item(13, 4).Value = ''
$ws.Cells.Item(13, 5).Value = ''
$ws.Cells.Item(13, 6).Value = ''
$ws.Cells.Item(13, 7).Value = 'X'
$ws.Cells.Item(13, 8).Value = ''
$ws.Cells.Item(13, 9).Value = ''
$ws.Cells.Item(13, 10).Value = ''

# Row 14
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = 'Trung tâm Ngoại ngữ - Tin học'
$ws.Cells.Item(14, 3).Value = '2023-2024'
$ws.Cells.Item(14, 4).Value = ''
$ws.Cells.Item(14, 5).Value = ''
$ws.Cells.Item(14, 6).Value = ''
$ws.Cells.Item(14, 7).Value = 'X'
$ws.Cells.Item(14, 8).Value = ''
$ws.Cells.Item(14, 9).Value = ''
$ws.Cells.Item(14, 10).Value = ''

# Row 15
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 'Trung tâm Ngoại ngữ - Tin học'
$ws.Cells.Item(15, 3).Value = '2023-2024'
$ws.Cells.Item(15, 4).Value = ''
$ws.Cells.Item(15, 5).Value = 'X'
$ws.Cells.Item(15, 6).Value = ''
$ws.Cells.Item(15, 7).Value = ''
$ws.Cells.Item(15, 8).Value = ''
$ws.Cells.Item(15, 9).Value = ''
$ws.Cells.Item(15, 10).Value = ''

# Row 16
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = 'Khoa Công nghệ thông tin - Truyền thông'
$ws.Cells.Item(16, 3).Value = '2024-2025'
$ws.Cells.Item(16, 4).Value = 'X'
$ws.Cells.Item(16, 5).Value = ''
$ws.Cells.Item(16, 6).Value = ''
$ws.Cells.Item(16, 7).Value = ''
$ws.Cells.Item(16, 8).Value = ''
$ws.Cells.Item(16, 9).Value = ''
$ws.Cells.Item(16, 10).Value = ''

# Row 17
$ws.Cells.Item(17, 1).Value = 16
$ws.Cells.Item(17, 2).Value = 'Trung tâm Ngoại ngữ - Tin học'
$ws.Cells.Item(17, 3).Value = '2024-2025'
$ws.Cells.Item(17, 4).Value = ''
$ws.Cells.Item(17, 5).Value = ''
$ws.Cells.Item(17, 6).Value = ''
$ws.Cells.Item(17, 7).Value = 'X'
$ws.Cells.Item(17, 8).Value = ''
$ws.Cells.Item(17, 9).Value = ''
$ws.Cells.Item(17, 10).Value = ''

# Row 18
$ws.Cells.Item(18, 1).Value = 17
$ws.Cells.Item(18, 2).Value = 'Trung tâm Ngoại ngữ - Tin học'
$ws.Cells.Item(18, 3).Value = '2024-2025'
$ws.Cells.Item(18, 4).Value = ''
$ws.Cells.Item(18, 5).Value = ''
$ws.Cells.Item(18, 6).Value = ''
$ws.Cells.Item(18, 7).Value = ''
$ws.Cells.Item(18, 8).Value = 'X'
$ws.Cells.Item(18, 9).Value = ''
$ws.Cells.Item(18, 10).Value = ''

# Row 19
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = 'Trung tâm Ngoại ngữ - Tin học'
$ws.Cells.Item(19, 3).Value = '2024-2025'
$ws.Cells.Item(19, 4).Value = ''
$ws.Cells.Item(19, 5).Value = ''
$ws.Cells.Item(19, 6).Value = ''
$ws.Cells.Item(19, 7).Value = ''
$ws.Cells.Item(19, 8).Value = ''
$ws.Cells.Item(19, 9).Value = 'X'
$ws.Cells.Item(19, 10).Value = ''

# Adjust column widths for B and C (best-effort, engine quantizes to pixel grid)
$ws.Columns.Item(2).ColumnWidth = 48.65
$ws.Columns.Item(3).ColumnWidth = 13.13

# Update selection to match the new used range
$ws.Range("A2:J19").Select()
